# Insert a new data row right before the current row 701, shifting the
# existing rows 701:742 down to 702:743 (dimension grows from D742 to D743).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(701).Insert()

# Fill in the newly inserted row with the new data point: 2026/01/26 (月), 16, 18
# Column A holds dates stored as plain text (matching the rest of the sheet,
# e.g. the "2026/01/26" entries in rows 698-700), so force text format before
# assigning the value to stop Excel from auto-converting the date-looking
# string into a real date serial number.
$a = $ws.Range("A701")
$a.NumberFormat = "@"
$a.Value = "2026/01/26"

$ws.Range("B701").Value = "月"
$ws.Range("C701").Value = 16
$ws.Range("D701").Value = 18
